$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update chartType values: "phyloTree" -> "phylogenetic_tree", "map" -> "geo"
# (order matters so the new shared strings land at the expected indices)
$ws.Range("A8").Value = "phylogenetic_tree"
$ws.Range("A6").Value = "geo"

# Move the active selection from A9 to A6
$ws.Range("A6").Select()
